$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1 headers: localize from Chinese to English
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Number"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Complex"
$ws.Range("E1").Value = "Ignored"
$ws.Range("F1").Value = "Empty"

# Row 2 template placeholders: keep the {.xxx} tokens, localize surrounding text
$ws.Range("A2").Value = "{.name}"
$ws.Range("B2").Value = "{.number}"
$ws.Range("C2").Value = "{.date}"
$ws.Range("D2").Value = "{.name} is {.number} years old this year"
$ws.Range("E2").Value = "\{.name\} ignored，{.name}"
$ws.Range("F2").Value = "Empty{.empty}"

# Move the active selection to E2 to match the updated demo state
$ws.Range("E2").Select()
